$d = $word.ActiveDocument

$pairs = @(
    @("93×70=", "60×63="),
    @("77×50=", "87×18="),
    @("68×70=", "22×79="),
    @("62×93=", "35×99="),
    @("49×97=", "65×50="),
    @("71×35=", "95×25="),
    @("83×23=", "58×83="),
    @("51×87=", "18×50="),
    @("99×53=", "52×53="),
    @("17×44=", "15×62="),
    @("33×77=", "51×84="),
    @("69×73=", "61×54="),
    @("37×73=", "24×98="),
    @("46×43=", "69×81="),
    @("50×76=", "96×11="),
    @("54×80=", "84×60="),
    @("46×67=", "19×51="),
    @("15×76=", "94×26="),
    @("50×16=", "24×61="),
    @("63×20=", "73×54="),
    @("31×20=", "26×60="),
    @("53×71=", "96×95="),
    @("58×82=", "29×68="),
    @("13×57=", "15×31="),
    @("86×56=", "89×55=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
